$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.662.93"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.491.97"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.76"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.97"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.70"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "67.566.34"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "2.514.60"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.55"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.23"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.67"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.77"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.15"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "2.622.58"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "0.0₃0905"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "509.83"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.84"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.36"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.329"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.02"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("E51").Value = "  +0.54%  "
